$wb = $excel.ActiveWorkbook

# --- Sheet 1: pages_with_related_resources ---
$ws1 = $wb.Worksheets.Item(1)

# Row 2 stays "about-cancer/coping/feelings / Article / English" but gets Wrap Text applied
$ws1.Range("A2").WrapText = $true
$ws1.Rows.Item(2).RowHeight = 17

# Row 3 - replace node/21 related content with the new vitamin-d blog page row
$ws1.Range("A3").Value = "about-cancer/coping/feelings/relaxation/vitamin-d-supplement-cancer-prevention"
$ws1.Range("B3").Value = "Blog Page"
$ws1.Range("C3").Value = "English"

# Row 4 - new Spanish blog page row
$ws1.Range("A4").Value = "espanol/about-cancer/coping/feelings/relaxation/vitamina-d-complemento-cancer-prevencion"
$ws1.Range("B4").Value = "Blog Page"
$ws1.Range("C4").Value = "Spanish"

# Row 5 - new press release row
$ws1.Range("A5").Value = "news-events/press-releases/2018/oropharyngeal-hpv-cisplatin"
$ws1.Range("B5").Value = "Press Release"
$ws1.Range("C5").Value = "English"

# Column widths
$ws1.Columns.Item(1).ColumnWidth = 73.83333333333334
$ws1.Columns.Item(2).ColumnWidth = 22.666666666666668

# --- Sheet 2: pages_without_related_resources ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 - replace node/21 with the Spanish "sobrellevar" article
$ws2.Range("A2").Value = "espanol/cancer/sobrellevar/sentimientos/hoja-informativa-estres"
$ws2.Range("B2").Value = "Article"
$ws2.Range("C2").Value = "Spanish"

# Row 3 - new blog row (alignment style touched but left at its default, matching
# the source workbook's empty "applyAlignment" cell style)
$ws2.Range("A3").Value = "news-events/cancer-currents-blog/2019/human-tumor-atlas-network-cancer-maps"
$ws2.Range("B3").Value = "Blog"
$ws2.Range("C3").Value = "English"
$ws2.Range("A3").WrapText = $true
$ws2.Range("A3").WrapText = $false

# Row 4 - new press release row (same touched-but-default alignment style)
$ws2.Range("A4").Value = "sharpless-nci-director"
$ws2.Range("B4").Value = "Press Release"
$ws2.Range("C4").Value = "English"
$ws2.Range("A4").WrapText = $true
$ws2.Range("A4").WrapText = $false

# Column width
$ws2.Columns.Item(1).ColumnWidth = 93.83333333333334

# --- Selections ---
$ws2.Range("A16").Select()
$ws1.Range("A12").Select()
